# Applies the "3e version avec organisation fichiers" edit to
# StructureDefinition-CompetenceMetier.xlsx:
#  - Metadata sheet: bump Date, change Base Definition URL
#  - Elements sheet: rename the "competenceMetier" element row to
#    "typeSavoirFaire", and add three new element rows (dateReconnaissance,
#    dateAbandon, competenceMetier) that describe the new SavoirFaire-based
#    structure.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-21T11:52:46+00:00"
$meta.Range("B18").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/SavoirFaire"

# ---------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Row 3 used to describe "CompetenceMetier.competenceMetier"; it is now
# repurposed to describe "CompetenceMetier.typeSavoirFaire".
$ws.Range("A3").Value = "CompetenceMetier.typeSavoirFaire"
$ws.Range("B3").Value = "CompetenceMetier.typeSavoirFaire"
$shortTypeSavoirFaire = " Le type de savoir-faire (qualifications/autres attributions) d" + [char]0x00E9 + "signe par exemple:** une sp" + [char]0x00E9 + "cialit" + [char]0x00E9 + " ordinale (S);** une comp" + [char]0x00E9 + "tence (C);** etc."
$ws.Range("L3").Value = $shortTypeSavoirFaire
$ws.Range("M3").Value = $shortTypeSavoirFaire
$ws.Range("Z3").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R04-TypeSavoirFaire?vs"
$ws.Range("AF3").Value = "SavoirFaire.typeSavoirFaire"

# Prepare rows 4-6 with the same formatting (borders / wrap / style) as
# row 3 before filling in their values.
$ws.Range("A3:AJ3").Copy()
$ws.Range("A4:AJ6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 4: CompetenceMetier.dateReconnaissance
$ws.Range("A4").Value = "CompetenceMetier.dateReconnaissance"
$ws.Range("B4").Value = "CompetenceMetier.dateReconnaissance"
$ws.Range("F4").Value = "0"
$ws.Range("G4").Value = "1"
$ws.Range("K4").Value = "date" + [char]10
$shortDateReconnaissance = " Date " + [char]0x00E0 + " laquelle, l" + [char]0x2019 + "organisme donnant l" + [char]0x2019 + "autorisation d" + [char]0x2019 + "exercer une qualification a reconnu cette qualification ou date " + [char]0x00E0 + " laquelle l'attribution a " + [char]0x00E9 + "t" + [char]0x00E9 + " donn" + [char]0x00E9 + "e au professionnel."
$ws.Range("L4").Value = $shortDateReconnaissance
$ws.Range("M4").Value = $shortDateReconnaissance
$ws.Range("X4").Value = ""
$ws.Range("Z4").Value = ""
$ws.Range("AF4").Value = "SavoirFaire.dateReconnaissance"
$ws.Range("AG4").Value = "0"
$ws.Range("AH4").Value = "1"

# Row 5: CompetenceMetier.dateAbandon
$ws.Range("A5").Value = "CompetenceMetier.dateAbandon"
$ws.Range("B5").Value = "CompetenceMetier.dateAbandon"
$ws.Range("F5").Value = "0"
$ws.Range("G5").Value = "1"
$ws.Range("K5").Value = "date" + [char]10
$shortDateAbandon = " Date " + [char]0x00E0 + " laquelle le professionnel a d" + [char]0x00E9 + "clar" + [char]0x00E9 + " renoncer " + [char]0x00E0 + " l" + [char]0x2019 + "exercice d" + [char]0x2019 + "un savoir-faire ou date " + [char]0x00E0 + " laquelle il ne souhaite plus le faire appara" + [char]0x00EE + "tre."
$ws.Range("L5").Value = $shortDateAbandon
$ws.Range("M5").Value = $shortDateAbandon
$ws.Range("X5").Value = ""
$ws.Range("Z5").Value = ""
$ws.Range("AF5").Value = "SavoirFaire.dateAbandon"
$ws.Range("AG5").Value = "0"
$ws.Range("AH5").Value = "1"

# Row 6: CompetenceMetier.competenceMetier (the original element, moved
# here and re-pointed at the new ValueSet).
$ws.Range("A6").Value = "CompetenceMetier.competenceMetier"
$ws.Range("B6").Value = "CompetenceMetier.competenceMetier"
$ws.Range("F6").Value = "0"
$ws.Range("G6").Value = "1"
$ws.Range("K6").Value = "Coding" + [char]10
$shortCompetenceMetier = " Comp" + [char]0x00E9 + "tence m" + [char]0x00E9 + "tier acquise par le professionnel"
$ws.Range("L6").Value = $shortCompetenceMetier
$ws.Range("M6").Value = $shortCompetenceMetier
$ws.Range("X6").Value = "preferred"
$ws.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/ValueSet/competenceMetier-vs"
$ws.Range("AF6").Value = "CompetenceMetier.competenceMetier"
$ws.Range("AG6").Value = "0"
$ws.Range("AH6").Value = "1"

# Keep the bestFit columns looking right after the new, wider content.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(26).AutoFit()
